$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("QuantitativeMetrics")

# Update B11 (BLEU score)
$ws.Range("B11").Value = 0.1312184895043395

# Update B12 (Code BLEU) and its related note in C12
$ws.Range("B12").Value = 0.3578968655708283
$ws.Range("C12").Value = "{'codebleu': 0.35789686557082834, 'ngram_match_score': 0.13121848950433948, 'weighted_ngram_match_score': 0.14399790088917605, 'syntax_match_score': 0.6237623762376238, 'dataflow_match_score': 0.532608695652174}"

# Update B13 (Embeddings and Cosine similarity)
$ws.Range("B13").Value = 0.9082444610469762
